$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3310
$ws.Range("F5").Value = 1357
$ws.Range("F8").Value = 388
$ws.Range("F10").Value = 52
$ws.Range("F11").Value = 8482
$ws.Range("F12").Value = 8482
$ws.Range("F13").Value = 461
$ws.Range("F16").Value = 92
$ws.Range("F17").Value = 311
$ws.Range("F19").Value = 76
$ws.Range("F21").Value = 348
$ws.Range("F22").Value = 10666
$ws.Range("F23").Value = 10666
$ws.Range("F27").Value = 142
$ws.Range("F31").Value = 136
$ws.Range("F32").Value = 2666
$ws.Range("F35").Value = 2073
$ws.Range("F36").Value = 35
$ws.Range("F39").Value = 886
$ws.Range("F40").Value = 4065
$ws.Range("F41").Value = 126
$ws.Range("F43").Value = 2572
$ws.Range("F44").Value = 3016
$ws.Range("F45").Value = 1229
$ws.Range("F48").Value = 329
$ws.Range("F49").Value = 290
$ws.Range("F51").Value = 110

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 75

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 3310
$ws.Range("F8").Value = 1358
$ws.Range("F11").Value = 388
$ws.Range("F16").Value = 8482
$ws.Range("F17").Value = 461
$ws.Range("F20").Value = 92
$ws.Range("F21").Value = 311
$ws.Range("F23").Value = 76
$ws.Range("F25").Value = 10666
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 75
$ws.Range("F29").Value = 142
$ws.Range("F34").Value = 136
$ws.Range("F35").Value = 2666
$ws.Range("F38").Value = 2073
$ws.Range("F39").Value = 35
$ws.Range("F42").Value = 886
$ws.Range("F44").Value = 126
$ws.Range("F45").Value = 3016
$ws.Range("F47").Value = 1229
$ws.Range("F48").Value = 329
$ws.Range("F49").Value = 290
$ws.Range("F51").Value = 110
